$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = 1000
$ws.Range("C3").Value = 1000
$ws.Range("C5").Value = 500
